$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Row 8: update pooh (H8) and reb (J8)
$ws.Range("H8").Value = 6
$ws.Range("J8").Value = 4

# Rows 68 and 69: swap player stat lines (columns D through O)
# Row 68 currently: Shawn Phillips Jr. / MIZ / AUB@MIZ / Final / 6 / 2 / 3 / 0 / 0 / 2 / 1 / 19
# Row 69 currently: Justin Abson / UGA / MISS@UGA / Final/OT / 5 / 4 / 0 / 0 / 1 / 0 / 0 / 7
# After the edit, row 68 should hold Justin Abson's line and row 69 should hold Shawn Phillips Jr.'s line.

$ws.Range("D68").Value = "Justin Abson"
$ws.Range("E68").Value = "UGA"
$ws.Range("F68").Value = "MISS@UGA"
$ws.Range("G68").Value = "Final/OT"
$ws.Range("H68").Value = 5
$ws.Range("I68").Value = 4
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 7

$ws.Range("D69").Value = "Shawn Phillips Jr."
$ws.Range("E69").Value = "MIZ"
$ws.Range("F69").Value = "AUB@MIZ"
$ws.Range("G69").Value = "Final"
$ws.Range("H69").Value = 5
$ws.Range("I69").Value = 2
$ws.Range("J69").Value = 3
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 2
$ws.Range("N69").Value = 1
$ws.Range("O69").Value = 19
